$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.279.12"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.589.78"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "212.17"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.812.87"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.632.98"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "64.44"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "26.280.22"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D20").Value = "213.86"
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Value = "8.98"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  -1.83%  "
$ws.Range("D25").Value = "'145.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "15.19"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "0.0501"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").Value = "1.341.23"
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "0.595"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "1.03"
$ws.Range("E40").Value = "  +15.15%  "
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  +3.81%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("D45").Value = "61.83"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "1.724.27"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "88.07"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  -4.45%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.0979"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0502"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("E51").Value = "  -0.37%  "
